# Natmi following Dr Hou advice
# Adds a new "ECs" sending-cluster block (and recomputed stats for the
# existing FAPs/M2/sCs sending clusters) to the Lamb2-Rpsa LR-pairs sheet.
# Rows 2-17 become a full 4 (sending cluster) x 4 (target cluster) matrix:
# ECs, FAPs, M2, sCs -> ECs, FAPs, M2, sCs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("ECs","Lamb2","Rpsa","ECs",2,0.6666666666666666,23.17518866666667,69.525566,0.2781252059919919,0.2781252059919919,3,1,69.72416666666666,209.1725,0.122848749818507,0.122848749818507,1615.870717126111,14542.836454135,0.03416733384913094,0.03416733384913094)
    ,@("ECs","Lamb2","Rpsa","FAPs",2,0.6666666666666666,23.17518866666667,69.525566,0.2781252059919919,0.2781252059919919,3,1,142.7363333333333,428.209,0.2514907089174394,0.2514907089174394,3307.941454588222,29771.473091294,0.06994590522273489,0.06994590522273489)
    ,@("ECs","Lamb2","Rpsa","M2",2,0.6666666666666666,23.17518866666667,69.525566,0.2781252059919919,0.2781252059919919,3,1,158.0830736666667,474.249221,0.2785305138211326,0.2785305138211326,3663.605057231565,32972.44551508409,0.07746635653155785,0.07746635653155785)
    ,@("ECs","Lamb2","Rpsa","sCs",2,0.6666666666666666,23.17518866666667,69.525566,0.2781252059919919,0.2781252059919919,3,1,197.0174863333333,591.052459,0.347130027442921,0.347130027442921,4565.917416407421,41093.25674766679,0.09654561038856818,0.09654561038856818)
    ,@("FAPs","Lamb2","Rpsa","ECs",3,1,38.67778666666666,116.03336,0.4641717285975485,0.4641717285975485,3,1,69.72416666666666,209.1725,0.122848749818507,0.122848749818507,2696.776443844444,24270.9879946,0.05702291655930419,0.05702291655930419)
    ,@("FAPs","Lamb2","Rpsa","FAPs",3,1,38.67778666666666,116.03336,0.4641717285975485,0.4641717285975485,3,1,142.7363333333333,428.209,0.2514907089174394,0.2514907089174394,5520.725450248889,49686.52905224,0.1167348770844307,0.1167348770844307)
    ,@("FAPs","Lamb2","Rpsa","M2",3,1,38.67778666666666,116.03336,0.4641717285975485,0.4641717285975485,3,1,158.0830736666667,474.249221,0.2785305138211326,0.2785305138211326,6114.303398890284,55028.73059001256,0.1292859900675185,0.1292859900675185)
    ,@("FAPs","Lamb2","Rpsa","sCs",3,1,38.67778666666666,116.03336,0.4641717285975485,0.4641717285975485,3,1,197.0174863333333,591.052459,0.347130027442921,0.347130027442921,7620.200306003581,68581.80275403224,0.1611279448862951,0.1611279448862951)
    ,@("M2","Lamb2","Rpsa","ECs",3,1,0.1530703333333333,0.459211,0.001836995530087286,0.001836995530087286,3,1,69.72416666666666,209.1725,0.122848749818507,0.122848749818507,10.67270143305556,96.05431289749998,0.0002256726042934087,0.0002256726042934087)
    ,@("M2","Lamb2","Rpsa","FAPs",3,1,0.1530703333333333,0.459211,0.001836995530087286,0.001836995530087286,3,1,142.7363333333333,428.209,0.2514907089174394,0.2514907089174394,21.84869812211111,196.638283099,0.0004619873081398188,0.0004619873081398188)
    ,@("M2","Lamb2","Rpsa","M2",3,1,0.1530703333333333,0.459211,0.001836995530087286,0.001836995530087286,3,1,158.0830736666667,474.249221,0.2785305138211326,0.2785305138211326,24.19782878051456,217.780459024631,0.0005116593088823356,0.0005116593088823356)
    ,@("M2","Lamb2","Rpsa","sCs",3,1,0.1530703333333333,0.459211,0.001836995530087286,0.001836995530087286,3,1,197.0174863333333,591.052459,0.347130027442921,0.347130027442921,30.15753230553878,271.417790749849,0.0006376763087717227,0.0006376763087717227)
    ,@("sCs","Lamb2","Rpsa","ECs",3,1,21.32041366666667,63.961241,0.2558660698803723,0.2558660698803723,3,1,69.72416666666666,209.1725,0.122848749818507,0.122848749818507,1486.548075896944,13378.9326830725,0.0314328268057785,0.0314328268057785)
    ,@("sCs","Lamb2","Rpsa","FAPs",3,1,21.32041366666667,63.961241,0.2558660698803723,0.2558660698803723,3,1,142.7363333333333,428.209,0.2514907089174394,0.2514907089174394,3043.197671929889,27388.779047369,0.06434793930213391,0.06434793930213391)
    ,@("sCs","Lamb2","Rpsa","M2",3,1,21.32041366666667,63.961241,0.2558660698803723,0.2558660698803723,3,1,158.0830736666667,474.249221,0.2785305138211326,0.2785305138211326,3370.396524271474,30333.56871844326,0.07126650791317393,0.07126650791317393)
    ,@("sCs","Lamb2","Rpsa","sCs",3,1,21.32041366666667,63.961241,0.2558660698803723,0.2558660698803723,3,1,197.0174863333333,591.052459,0.347130027442921,0.347130027442921,4200.494308193513,37804.44877374162,0.08881879585928598,0.08881879585928598)
)

$startRow = 2
$r = $startRow
foreach ($row in $data) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
